# Updates the "想去人数" (interest count, column F) values on the
# 展览 (Exhibition), 本地生活 (Local Life) and 全部类型 (All Types) sheets
# to reflect refreshed counts from a later data pull (gh-pages output
# regenerated at commit 456a3b4). The 演出 (Performance) sheet is
# unchanged.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    3  = 672
    7  = 649
    8  = 2210
    11 = 233
    14 = 1064
    15 = 427
    19 = 4403
    21 = 3362
    23 = 60
    25 = 3281
    26 = 4888
    28 = 970
    30 = 3170
    31 = 337
    37 = 1388
    39 = 1315
    42 = 787
    43 = 491
    44 = 51
    45 = 281
    46 = 57
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# --- Sheet: 本地生活 (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 2083

# --- Sheet: 全部类型 (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 2083
    3  = 672
    7  = 649
    9  = 2210
    12 = 233
    16 = 1064
    17 = 427
    19 = 4404
    23 = 3362
    24 = 3281
    25 = 4888
    27 = 970
    28 = 3170
    29 = 337
    35 = 1388
    37 = 1315
    40 = 491
    42 = 51
    44 = 281
    46 = 57
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
